$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 175, shifting existing rows 175:293 down to 176:294
$ws.Rows("175:175").Insert()

# Populate the newly inserted row 175 with the new weekly price entry.
# Columns that stay the same as the (now shifted) surrounding data:
$ws.Range("A175").Value = 3
$ws.Range("B175").Value = "Femacal de La Calera"
$ws.Range("C175").Value = "Coquimbo"
$ws.Range("D175").Value = 44603
$ws.Range("E175").Value = 5
$ws.Range("F175").Value = 100112043
$ws.Range("G175").Value = "Pepino ensalada"
$ws.Range("H175").Value = "Sin especificar"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 85
$ws.Range("K175").Value = 13000
$ws.Range("L175").Value = 14000
$ws.Range("M175").Value = 13588
$ws.Range("N175").Value = '$/caja 70 unidades'
$ws.Range("O175").Value = "Limache"
$ws.Range("P175").Value = 194
$ws.Range("Q175").Value = 70
$ws.Range("R175").Value = "Hortaliza"
